# Append 19 new data rows (140-158) to the "dataset" sheet, replicating the
# tensorflow/ranking repository entries that were appended when the file
# download API was implemented.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 139 (id 138) carries the "normal" bordered/bold/centered id-column
# style used by every data row from row 4 onward. Copy that formatting down
# onto the new id cells (column A) of the rows we are about to add.
$ws.Range("A139").Copy()
$ws.Range("A140:A158").PasteSpecial(-4122)

# The source data stores the date and the 0/1 flag columns as plain text
# (not real dates/numbers), e.g. "12/03/2018", "0", "1". Pre-format the new
# cells as Text so Excel does not silently coerce them into a date serial
# number or a numeric value when we assign the literal strings below.
$ws.Range("E140:O158").NumberFormat = "@"

$url = "https://github.com/tensorflow/ranking"
$repo = "ranking"
$author = "tensorflow"
$startDate = "12/03/2018"

for ($i = 0; $i -lt 19; $i++) {
    $r = 140 + $i
    $id = 139 + $i

    $ws.Cells.Item($r, 1).Value = $id
    $ws.Cells.Item($r, 2).Value = $url
    $ws.Cells.Item($r, 3).Value = $repo
    $ws.Cells.Item($r, 4).Value = $author
    $ws.Cells.Item($r, 5).Value = $startDate

    # OSE, BCE, PDE
    $ws.Cells.Item($r, 6).Value = "0"
    $ws.Cells.Item($r, 7).Value = "0"
    $ws.Cells.Item($r, 8).Value = "0"
    # SV
    $ws.Cells.Item($r, 9).Value = "1"
    # OS: 0 for the first two new rows (140-141), 1 afterwards (142-158)
    if ($r -le 141) {
        $ws.Cells.Item($r, 10).Value = "0"
    } else {
        $ws.Cells.Item($r, 10).Value = "1"
    }
    # SD
    $ws.Cells.Item($r, 11).Value = "1"
    # RS, TFS
    $ws.Cells.Item($r, 12).Value = "0"
    $ws.Cells.Item($r, 13).Value = "0"
    # UI
    $ws.Cells.Item($r, 14).Value = "1"
    # TC
    $ws.Cells.Item($r, 15).Value = "0"
}
